# Insert a new data row before current row 465 (shifting existing rows 465-540 down to 466-541),
# populate the new row with a copy of the (now shifted) old row 465 data, then set the
# new row's date (D) and volume (J) values to reflect the newly added record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 465; rows 465:540 move to 466:541
$ws.Rows("465:465").Insert()

# The data that used to live in row 465 is now in row 466. Copy it into the new
# blank row 465 so all the non-changing columns (A,B,C,E,F,G,H,I,K,L,M,N,O,P,Q,R)
# are populated correctly.
$ws.Rows("466:466").Copy()
$ws.Rows("465:465").PasteSpecial()
$excel.CutCopyMode = 0

# Apply the new values for this inserted record.
$ws.Range("D465").Value = 45218
$ws.Range("J465").Value = 45
